# Apply the Feb-15-2023 crypto symbol-list refresh described by the commit.
#
# Source sheet stores every data cell as literal text (t="inlineStr"), including
# numeric-looking values such as prices, percentages and the refresh-minute column.
# A plain `.Value = "296.85"` assignment would let Excel re-interpret that as a
# *number*, changing the cell's stored type relative to the source workbook.
# Prefixing with an apostrophe (the same trick Excel's UI uses for "text that
# looks numeric") keeps those cells as text, matching the original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = "'"

# Row 2
$ws.Cells.Item(2, 4).Value = $apos + "296.85"
$ws.Cells.Item(2, 5).Value = $apos + "1.68%"
$ws.Cells.Item(2, 7).Value = $apos + "6"
# Row 3
$ws.Cells.Item(3, 4).Value = $apos + "41.92"
$ws.Cells.Item(3, 5).Value = $apos + "3.87%"
$ws.Cells.Item(3, 7).Value = $apos + "6"
# Row 4
$ws.Cells.Item(4, 4).Value = $apos + "5.018"
$ws.Cells.Item(4, 5).Value = $apos + "-0.26%"
$ws.Cells.Item(4, 7).Value = $apos + "6"
# Row 5
$ws.Cells.Item(5, 4).Value = $apos + "0.07520"
$ws.Cells.Item(5, 5).Value = $apos + "2.62%"
$ws.Cells.Item(5, 7).Value = $apos + "6"
# Row 6
$ws.Cells.Item(6, 2).Value = "FTXToken"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(6, 4).Value = $apos + "1.582"
$ws.Cells.Item(6, 5).Value = $apos + "3.01%"
$ws.Cells.Item(6, 7).Value = $apos + "6"
# Row 7
$ws.Cells.Item(7, 2).Value = "MXToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(7, 4).Value = $apos + "0.9256"
$ws.Cells.Item(7, 5).Value = $apos + "-0.15%"
$ws.Cells.Item(7, 7).Value = $apos + "6"
# Row 8
$ws.Cells.Item(8, 2).Value = "BTSEToken"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(8, 4).Value = $apos + "2.401"
$ws.Cells.Item(8, 5).Value = $apos + "2.26%"
$ws.Cells.Item(8, 7).Value = $apos + "6"
# Row 9
$ws.Cells.Item(9, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(9, 4).Value = $apos + "0.1198"
$ws.Cells.Item(9, 5).Value = $apos + "2.19%"
$ws.Cells.Item(9, 7).Value = $apos + "6"
# Row 10
$ws.Cells.Item(10, 2).Value = "WazirX"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10, 4).Value = $apos + "0.1825"
$ws.Cells.Item(10, 5).Value = $apos + "4.76%"
$ws.Cells.Item(10, 7).Value = $apos + "6"
# Row 11
$ws.Cells.Item(11, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(11, 4).Value = $apos + "0.08959"
$ws.Cells.Item(11, 5).Value = $apos + "3.21%"
$ws.Cells.Item(11, 7).Value = $apos + "6"
# Row 12
$ws.Cells.Item(12, 2).Value = "BitrueCoin"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(12, 4).Value = $apos + "0.04068"
$ws.Cells.Item(12, 5).Value = $apos + "-6.16%"
$ws.Cells.Item(12, 7).Value = $apos + "6"
# Row 13
$ws.Cells.Item(13, 2).Value = "BitMartToken"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(13, 4).Value = $apos + "0.1050"
$ws.Cells.Item(13, 5).Value = $apos + "-0.58%"
$ws.Cells.Item(13, 7).Value = $apos + "6"
# Row 14
$ws.Cells.Item(14, 2).Value = "BitForexToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(14, 4).Value = $apos + "0.001293"
$ws.Cells.Item(14, 5).Value = $apos + "1.95%"
$ws.Cells.Item(14, 7).Value = $apos + "6"
# Row 15
$ws.Cells.Item(15, 2).Value = "TigerCash"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(15, 4).Value = $apos + "0.005908"
$ws.Cells.Item(15, 5).Value = $apos + "-2.05%"
$ws.Cells.Item(15, 7).Value = $apos + "6"
# Row 16
$ws.Cells.Item(16, 2).Value = "LEO"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(16, 4).Value = $apos + "3.356"
$ws.Cells.Item(16, 5).Value = $apos + "0.55%"
$ws.Cells.Item(16, 7).Value = $apos + "6"
# Row 17
$ws.Cells.Item(17, 2).Value = "GateToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(17, 4).Value = $apos + "4.382"
$ws.Cells.Item(17, 5).Value = $apos + "2.01%"
$ws.Cells.Item(17, 7).Value = $apos + "6"
# Row 18
$ws.Cells.Item(18, 7).Value = $apos + "6"
# Row 19
$ws.Cells.Item(19, 4).Value = $apos + "8.065"
$ws.Cells.Item(19, 5).Value = $apos + "1.15%"
$ws.Cells.Item(19, 7).Value = $apos + "6"
# Row 20
$ws.Cells.Item(20, 5).Value = $apos + "-3.53%"
$ws.Cells.Item(20, 7).Value = $apos + "6"
# Row 21
$ws.Cells.Item(21, 5).Value = $apos + "13.08%"
$ws.Cells.Item(21, 7).Value = $apos + "6"
# Row 22
$ws.Cells.Item(22, 4).Value = $apos + "0.04090"
$ws.Cells.Item(22, 5).Value = $apos + "3.86%"
$ws.Cells.Item(22, 7).Value = $apos + "6"
# Row 23
$ws.Cells.Item(23, 5).Value = $apos + "0.32%"
$ws.Cells.Item(23, 7).Value = $apos + "6"
# Row 24
$ws.Cells.Item(24, 4).Value = $apos + "0.003939"
$ws.Cells.Item(24, 5).Value = $apos + "4.22%"
$ws.Cells.Item(24, 7).Value = $apos + "6"
# Row 25
$ws.Cells.Item(25, 5).Value = $apos + "-3.94%"
$ws.Cells.Item(25, 7).Value = $apos + "6"
# Row 26
$ws.Cells.Item(26, 7).Value = $apos + "6"
# Row 27
$ws.Cells.Item(27, 7).Value = $apos + "6"
# Row 28
$ws.Cells.Item(28, 7).Value = $apos + "6"
# Row 29
$ws.Cells.Item(29, 7).Value = $apos + "6"
# Row 30
$ws.Cells.Item(30, 7).Value = $apos + "6"
# Row 31
$ws.Cells.Item(31, 7).Value = $apos + "6"
# Row 32
$ws.Cells.Item(32, 7).Value = $apos + "6"
# Row 33
$ws.Cells.Item(33, 7).Value = $apos + "6"
# Row 34
$ws.Cells.Item(34, 7).Value = $apos + "6"
# Row 35
$ws.Cells.Item(35, 7).Value = $apos + "6"
# Row 36
$ws.Cells.Item(36, 7).Value = $apos + "6"
# Row 37
$ws.Cells.Item(37, 7).Value = $apos + "6"
# Row 38
$ws.Cells.Item(38, 4).Value = $apos + "0.02408"
$ws.Cells.Item(38, 5).Value = $apos + "5.80%"
$ws.Cells.Item(38, 7).Value = $apos + "6"
# Row 39
$ws.Cells.Item(39, 4).Value = $apos + "0.05206"
$ws.Cells.Item(39, 5).Value = $apos + "4.34%"
$ws.Cells.Item(39, 7).Value = $apos + "6"
# Row 40
$ws.Cells.Item(40, 4).Value = $apos + "0.006306"
$ws.Cells.Item(40, 5).Value = $apos + "6.66%"
$ws.Cells.Item(40, 7).Value = $apos + "6"
# Row 41
$ws.Cells.Item(41, 4).Value = $apos + "0.007816"
$ws.Cells.Item(41, 5).Value = $apos + "1.36%"
$ws.Cells.Item(41, 7).Value = $apos + "6"
# Row 42
$ws.Cells.Item(42, 4).Value = $apos + "0.1326"
$ws.Cells.Item(42, 5).Value = $apos + "3.24%"
$ws.Cells.Item(42, 7).Value = $apos + "6"
# Row 43
$ws.Cells.Item(43, 4).Value = $apos + "0.007404"
$ws.Cells.Item(43, 5).Value = $apos + "0.60%"
$ws.Cells.Item(43, 7).Value = $apos + "6"
# Row 44
$ws.Cells.Item(44, 4).Value = $apos + "0.007834"
$ws.Cells.Item(44, 5).Value = $apos + "-5.48%"
$ws.Cells.Item(44, 7).Value = $apos + "6"
# Row 45
$ws.Cells.Item(45, 4).Value = $apos + "0.2960"
$ws.Cells.Item(45, 5).Value = $apos + "1.51%"
$ws.Cells.Item(45, 7).Value = $apos + "6"
# Row 46
$ws.Cells.Item(46, 4).Value = $apos + "0.00006600"
$ws.Cells.Item(46, 5).Value = $apos + "5.07%"
$ws.Cells.Item(46, 7).Value = $apos + "6"
# Row 47
$ws.Cells.Item(47, 5).Value = $apos + "-0.03%"
$ws.Cells.Item(47, 7).Value = $apos + "6"
# Row 48
$ws.Cells.Item(48, 4).Value = $apos + "0.04490"
$ws.Cells.Item(48, 5).Value = $apos + "15.91%"
$ws.Cells.Item(48, 7).Value = $apos + "6"
# Row 49
$ws.Cells.Item(49, 4).Value = $apos + "0.004204"
$ws.Cells.Item(49, 5).Value = $apos + "0.05%"
$ws.Cells.Item(49, 7).Value = $apos + "6"
# Row 50
$ws.Cells.Item(50, 5).Value = $apos + "-0.03%"
$ws.Cells.Item(50, 7).Value = $apos + "6"
# Row 51
$ws.Cells.Item(51, 5).Value = $apos + "-0.03%"
$ws.Cells.Item(51, 7).Value = $apos + "6"
